# "separate dept from affiliations"
#
# - "PI hours" keeps name/hours/percentage/dept but gains a new "app"
#   column holding the full (possibly multi-valued) affiliation list that
#   used to live in the "dept" column; "dept" itself is narrowed down to
#   the PI's single primary department.
# - "dept hours" is renamed to "department hours" and now tallies hours
#   per *primary* dept (one row per PI, same hours/percentage as "PI hours").
# - A new "unit(accumulative) hours" sheet is appended holding the old
#   "dept hours" breakdown (hours accumulated across every affiliation).

$wb = $excel.ActiveWorkbook

$piSheet = $wb.Worksheets.Item("PI hours")
$deptSheet = $wb.Worksheets.Item("dept hours")

# --- capture the old "dept hours" values before we overwrite them ------
$oldRow2Dept = $deptSheet.Cells.Item(2, 2).Text
$oldRow2Hours = $deptSheet.Cells.Item(2, 3).Value()
$oldRow2Pct = $deptSheet.Cells.Item(2, 4).Value()
$oldRow3Dept = $deptSheet.Cells.Item(3, 2).Text
$oldRow3Hours = $deptSheet.Cells.Item(3, 3).Value()
$oldRow3Pct = $deptSheet.Cells.Item(3, 4).Value()
$oldRow4Dept = $deptSheet.Cells.Item(4, 2).Text
$oldRow4Hours = $deptSheet.Cells.Item(4, 3).Value()
$oldRow4Pct = $deptSheet.Cells.Item(4, 4).Value()

# --- capture the old "PI hours" per-PI dept/affiliation values ---------
$oldPiRow2Dept = $piSheet.Cells.Item(2, 5).Text
$oldPiRow3Dept = $piSheet.Cells.Item(3, 5).Text

# 1) "PI hours": add the "app" column (old multi-valued affiliation list)
#    and shrink "dept" down to a single department per PI.
$piSheet.Cells.Item(1, 6).Value = "app"
$piSheet.Cells.Item(2, 6).Value = $oldPiRow2Dept
$piSheet.Cells.Item(3, 6).Value = $oldPiRow3Dept

$piSheet.Range("B1").Copy()
$piSheet.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$piSheet.Cells.Item(2, 5).Value = "ME"
$piSheet.Cells.Item(3, 5).Value = "AE"

# 2) Append the new "unit(accumulative) hours" sheet with the old
#    "dept hours" breakdown, then rename "dept hours" itself.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$unitSheet = $wb.Worksheets.Add($null, $lastSheet)
$unitSheet.Name = "unit(accumulative) hours"

$unitSheet.Cells.Item(1, 2).Value = "unit(accumulative)"
$unitSheet.Cells.Item(1, 3).Value = "hours"
$unitSheet.Cells.Item(1, 4).Value = "percentage"

$unitSheet.Cells.Item(2, 1).Value = 0
$unitSheet.Cells.Item(2, 2).Value = $oldRow2Dept
$unitSheet.Cells.Item(2, 3).Value = $oldRow2Hours
$unitSheet.Cells.Item(2, 4).Value = $oldRow2Pct

$unitSheet.Cells.Item(3, 1).Value = 1
$unitSheet.Cells.Item(3, 2).Value = $oldRow3Dept
$unitSheet.Cells.Item(3, 3).Value = $oldRow3Hours
$unitSheet.Cells.Item(3, 4).Value = $oldRow3Pct

$unitSheet.Cells.Item(4, 1).Value = 2
$unitSheet.Cells.Item(4, 2).Value = $oldRow4Dept
$unitSheet.Cells.Item(4, 3).Value = $oldRow4Hours
$unitSheet.Cells.Item(4, 4).Value = $oldRow4Pct

$deptSheet.Range("B1").Copy()
$unitSheet.Range("B1:D1").PasteSpecial(-4122)
$deptSheet.Range("A2:A3").Copy()
$unitSheet.Range("A2:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Overwrite "dept hours" with the new per-PI-primary-dept breakdown,
#    then rename it to "department hours".
$deptSheet.Cells.Item(2, 2).Value = "ME"
$deptSheet.Cells.Item(2, 3).Value = 61
$deptSheet.Cells.Item(2, 4).Value = 88.40579710144928

$deptSheet.Cells.Item(3, 2).Value = "AE"
$deptSheet.Cells.Item(3, 3).Value = 8
$deptSheet.Cells.Item(3, 4).Value = 11.59420289855072

$deptSheet.Rows.Item(4).Delete()

$deptSheet.Name = "department hours"

# --- restore "PI hours" as the active/selected sheet --------------------
$piSheet.Activate()
